# edit.ps1
# Applies the workbook changes:
#  1. Rename "Weekly Quantity" B1 header from "Requested quantity" -> "Weekly_PO_Qty"
#  2. Rename "Monthly Trend" B1 header from "Requested quantity" -> "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with ds / PO_Forecast /
#     yhat_lower / yhat_upper columns and 79 rows of forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity header rename ---------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend header rename -------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the PO Forecast worksheet -----------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used on the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in).
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

$data = @"
44948.99999999999;45;-636.6430257260394;723.5884720458748
44969.99999999999;68;-589.6831212151328;738.5850709908807
44976.99999999999;75;-571.2261679555558;749.4075569795202
44983.99999999999;83;-585.7281483827235;746.5668945988839
44990.99999999999;91;-590.8244827857113;711.0431148223964
44997.99999999999;98;-535.0624971734999;719.7131259354635
45004.99999999999;106;-582.2913772473448;765.8290540136693
45011.99999999999;114;-563.5980361864137;796.8814818601787
45025.99999999999;129;-519.9238741707062;794.1111648844374
45060.99999999999;167;-470.3648577389931;810.2292965940923
45067.99999999999;175;-457.4302483824827;878.8702674516313
45074.99999999999;183;-467.492073604029;920.1501272542042
45095.99999999999;206;-426.8853135975946;861.0512635832769
45102.99999999999;213;-423.5340256473335;932.6078899155523
45109.99999999999;221;-463.180862909428;833.650622235897
45116.99999999999;229;-438.6021760031375;943.1719890998939
45123.99999999999;236;-459.8927931888003;946.0198835177441
45130.99999999999;244;-450.1868071293737;895.0770928017482
45137.99999999999;252;-419.6362044522973;923.8616731309685
45151.99999999999;267;-392.4952252503711;936.8398693637459
45165.99999999999;282;-314.2758046507714;971.4090404001299
45186.99999999999;305;-414.8463921836206;980.1141417636125
45193.99999999999;313;-383.9588805973181;960.1961906087387
45200.99999999999;321;-398.1468806889637;968.3370098517097
45207.99999999999;328;-300.4593911293841;993.4159362035201
45214.99999999999;336;-362.4780832236549;1048.008735011729
45221.99999999999;344;-332.1092700455959;1012.49814606384
45228.99999999999;351;-320.9892437962333;1009.358407533439
45235.99999999999;359;-295.6048878414422;1008.231633925392
45242.99999999999;367;-269.669842235748;1069.443945965432
45249.99999999999;374;-340.1264480370506;1021.82033576792
45256.99999999999;382;-304.0004402957453;1089.401693723235
45263.99999999999;390;-240.0781717320419;1023.52862616819
45277.99999999999;405;-281.8584971199132;1108.731542600141
45298.99999999999;428;-253.3684232501689;1124.964345462293
45305.99999999999;436;-242.3478502934478;1124.895098191723
45312.99999999999;443;-213.6281624515166;1103.872430939524
45319.99999999999;451;-191.4488150168629;1101.976895008428
45326.99999999999;459;-214.8125270659654;1132.074338246837
45333.99999999999;466;-165.093640363066;1131.594977956026
45340.99999999999;474;-197.4055398982362;1187.198368008229
45347.99999999999;482;-219.6495232183699;1146.616839746782
45354.99999999999;489;-185.9169183449868;1184.946277202754
45361.99999999999;497;-152.7401649613681;1226.898912649489
45368.99999999999;505;-146.0054812207898;1180.29805851638
45375.99999999999;512;-182.4807037775768;1185.52366181859
45382.99999999999;520;-192.0056020767541;1150.998916764872
45389.99999999999;528;-165.5644801351283;1254.170157941543
45396.99999999999;535;-108.7118957132869;1242.443094895355
45403.99999999999;543;-97.53238305789634;1261.309173800735
45410.99999999999;551;-124.0646052630558;1159.556291789233
45417.99999999999;558;-124.3841071550571;1252.810036639105
45424.99999999999;566;-117.2885040736451;1230.700246051374
45438.99999999999;581;-56.79747608366263;1269.985373966387
45445.99999999999;589;-62.60516577239999;1272.590006835926
45459.99999999999;604;-34.66948586059365;1267.083664573039
45466.99999999999;612;-63.37930539624385;1320.864016987359
45473.99999999999;620;-76.98061451978518;1262.873764146057
45487.99999999999;635;-74.19932502584328;1315.965414320755
45494.99999999999;643;-88.90773978567238;1308.956890662789
45501.99999999999;650;28.10009147332133;1306.53237021419
45515.99999999999;666;-34.1288965844203;1368.902094325354
45529.99999999999;681;16.85225285360239;1331.891437803813
45536.99999999999;689;-5.238889653378166;1420.34478702922
45543.99999999999;696;55.08320180158032;1348.522272933445
45550.99999999999;704;72.78076293855977;1393.668649413338
45557.99999999999;712;23.89716169152943;1345.059331549625
45592.99999999999;750;116.6971162436673;1409.977726954412
45599.99999999999;758;83.69377742135788;1437.044458020383
45613.99999999999;773;103.8782006108768;1462.673875422163
45627.99999999999;789;107.2114460117069;1451.087922556911
45634.99999999999;796;78.78999802917812;1483.282991728968
45641.99999999999;804;155.495618739007;1493.980164497453
45648.99999999999;812;140.8027248788518;1506.963170226894
45655.99999999999;819;129.9111838359548;1459.623355153187
45662.99999999999;827;176.2134215983134;1528.89820855122
45669.99999999999;835;152.5257951743156;1475.073831143323
45676.99999999999;842;155.3744851381666;1549.484075368446
45683.99999999999;850;170.2533249644576;1503.866383645433
"@

$lines = $data -split "`r?`n"
$rowCount = $lines.Count

$values = New-Object 'object[,]' $rowCount,4
for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $lines[$i].Split(";")
    $values[$i,0] = [double]$parts[0]
    $values[$i,1] = [double]$parts[1]
    $values[$i,2] = [double]$parts[2]
    $values[$i,3] = [double]$parts[3]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1

# Write header text first
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Bulk-write the data block
$wsForecast.Range("A$startRow" + ":D$endRow").Value2 = $values

# --- Formatting: match the look of the other two sheets -----------------
# Header row formatting (bold, centered, bordered) - copy from an existing
# header cell so the exact same style definition is reused.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Date column (A) formatting - copy the existing date-style cell.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A$endRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the originally-active sheet/tab selected (sheet add/activate shifts focus).
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null

Write-Host "PO Forecast sheet added with $rowCount data rows."
